# Remove all the thumbnail picture shapes from the slide, keeping only
# the first picture (the QR code, "Picture 3").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Walk the shapes from the end back to just after the QR code picture
# (index 3) and delete each one. Deleting from the tail avoids index
# shifting issues.
for ($i = $s.Shapes.Count; $i -ge 4; $i--) {
    $s.Shapes.Item($i).Delete()
}
